# Update the "want to go" count (column F) on several rows across the
# 展览 / 演出 / 本地生活 / 全部类型 sheets, matching the refreshed scrape
# output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 41671
$ws1.Range("F5").Value = 9421
$ws1.Range("F7").Value = 832
$ws1.Range("F8").Value = 832
$ws1.Range("F17").Value = 719
$ws1.Range("F29").Value = 491
$ws1.Range("F30").Value = 511
$ws1.Range("F33").Value = 920
$ws1.Range("F36").Value = 87
$ws1.Range("F39").Value = 372
$ws1.Range("F40").Value = 1233
$ws1.Range("F46").Value = 12

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 329

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 512
$ws3.Range("F4").Value = 377

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 41671
$ws4.Range("F7").Value = 329
$ws4.Range("F9").Value = 9421
$ws4.Range("F10").Value = 832
$ws4.Range("F11").Value = 832
$ws4.Range("F19").Value = 719
$ws4.Range("F31").Value = 491
$ws4.Range("F34").Value = 511
$ws4.Range("F40").Value = 87
$ws4.Range("F43").Value = 372
